$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (Changed) date in column C for every existing
#    data row (2-321) from 2023-09-13 (45182) to 2023-09-15 (45184).
for ($r = 2; $r -le 321; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45184
}

# 2) Row 321 picks up an explicit row height in the new file; give it the
#    same height as the rest of the data rows so the attribute is written.
$ws.Rows.Item(321).RowHeight = 15

# 3) Append the new record as row 322.
$ws.Cells.Item(322, 1).Value = "A 42797-2023"
$ws.Cells.Item(322, 2).Value2 = 45182
$ws.Cells.Item(322, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(322, 3).Value2 = 45184
$ws.Cells.Item(322, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(322, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(322, 5).Value = "ANEBY"
$ws.Cells.Item(322, 7).Value2 = 2
$ws.Cells.Item(322, 8).Value2 = 0
$ws.Cells.Item(322, 9).Value2 = 0
$ws.Cells.Item(322, 10).Value2 = 0
$ws.Cells.Item(322, 11).Value2 = 0
$ws.Cells.Item(322, 12).Value2 = 0
$ws.Cells.Item(322, 13).Value2 = 0
$ws.Cells.Item(322, 14).Value2 = 0
$ws.Cells.Item(322, 15).Value2 = 0
$ws.Cells.Item(322, 16).Value2 = 0
$ws.Cells.Item(322, 17).Value2 = 0
$ws.Cells.Item(322, 18).WrapText = $true
